$p = $ppt.ActivePresentation
$null = $p.Slides.Item(1).NotesPage
